$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 67: remove the empty B67 placeholder cell (becomes fully empty / absent)
$ws.Range("B67").ClearContents()

# Row 68
$ws.Range("A68").Value = "05/01/2026 10:58:53"
$ws.Range("B68").Value = "05/01 10:51"
$ws.Range("C68").Value = "Metrópoles"
$ws.Range("D68").Value = "Moraes manda PF explicar queixa de Bolsonaro sobre ar-condicionado"
$ws.Range("E68").Value = "https://www.metropoles.com/brasil/moraes-manda-pf-explicar-queixa-de-bolsonaro-sobre-ar-condicionado"
$ws.Range("F68").Value = "stf"
$ws.Range("G68").Value = "Defesa acionou o STF e pediu providências para corrigir barulho contínuo no local onde Bolsonaro está preso, n"

# Row 69
$ws.Range("A69").Value = "05/01/2026 10:58:54"
$ws.Range("B69").Value = "05/01 10:44"
$ws.Range("C69").Value = "g1 > Política"
$ws.Range("D69").Value = "Após queixas de Bolsonaro, Moraes manda PF prestar informações sobre ruídos em ar-condicionado"
$ws.Range("E69").Value = "https://g1.globo.com/politica/noticia/2026/01/05/apos-reclamacao-de-bolsonaro-moraes-manda-pf-prestar-informacoes-sobre-ruidos-em-ar-condicionado.ghtml"
$ws.Range("F69").Value = "stf"
$ws.Range("G69").Value = "O ministro Alexandre de Moraes, do Supremo Tribunal Federal (STF), determinou à Polícia Federal que preste esclarecimentos sobre ruídos no sistema de ar-c"

# Row 70
$ws.Range("A70").Value = "05/01/2026 10:58:55"
$ws.Range("B70").Value = "05/01 10:44"
$ws.Range("C70").Value = "Metrópoles"
$ws.Range("D70").Value = "Focus: mercado projeta alta na inflação e Selic em 12,25% em 2026"
$ws.Range("E70").Value = "https://www.metropoles.com/brasil/focus-mercado-projeta-alta-na-inflacao-e-selic-em-1225-em-2026"

# Row 71
$ws.Range("A71").Value = "05/01/2026 10:58:56"
$ws.Range("B71").Value = "05/01 10:35"
$ws.Range("C71").Value = "g1 > Economia"
$ws.Range("D71").Value = "Instabilidade na Venezuela impulsiona alta do ouro e recorde da prata"
$ws.Range("E71").Value = "https://g1.globo.com/economia/noticia/2026/01/05/instabilidade-na-venezuela-impulsiona-alta-do-ouro-e-recorde-da-prata.ghtml"
$ws.Range("F71").Value = "juros"
$ws.Range("G71").Value = "set Management`nOuro e prata costumam se destacar nesses momentos, especialmente quando os &lt;b&gt;juros&lt;/b&gt; estão mais baixos. Outros metais, como o cobre, também reagem ao cenário, refletindo tant"

# Row 72
$ws.Range("A72").Value = "05/01/2026 10:58:57"
$ws.Range("B72").Value = "05/01 10:26"
$ws.Range("C72").Value = "g1 > Política"
$ws.Range("D72").Value = "Governo Lula vê desordem global no ataque à Venezuela; Planalto avalia que Trump faz negócios com a ofensiva e monitora eventual risco à eleição no Brasil"
$ws.Range("E72").Value = "https://g1.globo.com/politica/blog/andreia-sadi/post/2026/01/05/governo-lula-venezuela-desordem-trump-eua-eleicao-brasil.ghtml"
$ws.Range("F72").Value = "lula"
$ws.Range("G72").Value = "t8DpQcewmHw/getattachmentthumbnail.png`" /&gt;&lt;br /&gt;     Integrantes do governo do presidente Lula (PT) avaliam como grave e preocupante o episódio envolvendo a Venezuela por representar o"

# Row 73
$ws.Range("A73").Value = "05/01/2026 10:59:01"
$ws.Range("C73").Value = "VEJA"
$ws.Range("D73").Value = "Índice de Preços ao Consumidor acelera e fecha 2025 com inflação de 4%"
$ws.Range("E73").Value = "https://veja.abril.com.br/economia/indice-de-precos-ao-consumidor-acelera-e-fecha-2025-com-inflacao-de-4/"
